# Remove the two placeholder team-member rows (name + ID) from the
# "Submitted By" table. The cell's run is deleted entirely (not just
# blanked) so the resulting paragraph has no <w:r> child, matching a
# normal Word "select cell contents, press Delete" edit.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Clear-CellText($table, $row, $col) {
    $cell = $table.Cell($row, $col)
    $rng = $cell.Range
    $txt = $rng.Text
    # $rng.Text includes the trailing paragraph mark + end-of-cell mark,
    # so trim those two characters off to get the real content length.
    $len = $txt.Length - 2
    if ($len -gt 0) {
        $delRange = $d.Range($rng.Start, $rng.Start + $len)
        $delRange.Delete()
    }
}

# Row 3: "Anindo Mahmood" / "23-55004-3"
Clear-CellText $t 3 1
Clear-CellText $t 3 2

# Row 4: "Nishat Tasnim Ema" / "23-55025-3"
Clear-CellText $t 4 1
Clear-CellText $t 4 2
